$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.266.73'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4973'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3978'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.10000'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +27.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.28'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.477'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.88'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.68'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.002'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.375'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001144'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.53%  '
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06653'
$ws.Range("D19").ClearFormats()
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.36'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.067'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.400.23'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.13%  '
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.24'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.069.43'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.474'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.70'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.34'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  -4.05%  '
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.644'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.587'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06801'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.261'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02382'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2165'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.022'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6283'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9996'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5986'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.688'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.71'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.978'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.193'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06837'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.04%  '
